$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.847215666666667
$ws.Range("H2").Value = 11.541647
$ws.Range("I2").Value = 0.05478559966737641
$ws.Range("J2").Value = 0.05478559966737641
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.407369
$ws.Range("N2").Value = 7.222107
$ws.Range("O2").Value = 0.05267229306167105
$ws.Range("P2").Value = 0.05267229306167105
$ws.Range("Q2").Value = 9.261667732247668
$ws.Range("R2").Value = 83.35500959022902
$ws.Range("S2").Value = 0.002885683161239438
$ws.Range("T2").Value = 0.002885683161239438

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.847215666666667
$ws.Range("H3").Value = 11.541647
$ws.Range("I3").Value = 0.05478559966737641
$ws.Range("J3").Value = 0.05478559966737641
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.67754066666667
$ws.Range("N3").Value = 35.032622
$ws.Range("O3").Value = 0.2555000268900398
$ws.Range("P3").Value = 0.2555000268900398
$ws.Range("Q3").Value = 44.92601740093712
$ws.Range("R3").Value = 404.3341566084341
$ws.Range("S3").Value = 0.01399772218820163
$ws.Range("T3").Value = 0.01399772218820163

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.847215666666667
$ws.Range("H4").Value = 11.541647
$ws.Range("I4").Value = 0.05478559966737641
$ws.Range("J4").Value = 0.05478559966737641
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.130105
$ws.Range("N4").Value = 9.390315
$ws.Range("O4").Value = 0.06848547433891598
$ws.Range("P4").Value = 0.06848547433891598
$ws.Range("Q4").Value = 12.04218899431167
$ws.Range("R4").Value = 108.379700948805
$ws.Range("S4").Value = 0.003752017780162231
$ws.Range("T4").Value = 0.003752017780162231

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.847215666666667
$ws.Range("H5").Value = 11.541647
$ws.Range("I5").Value = 0.05478559966737641
$ws.Range("J5").Value = 0.05478559966737641
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.949797
$ws.Range("N5").Value = 8.849391
$ws.Range("O5").Value = 0.06454040575268606
$ws.Range("P5").Value = 0.06454040575268606
$ws.Range("Q5").Value = 11.34850523188634
$ws.Range("R5").Value = 102.136547086977
$ws.Range("S5").Value = 0.003535884831936696
$ws.Range("T5").Value = 0.003535884831936696

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.847215666666667
$ws.Range("H6").Value = 11.541647
$ws.Range("I6").Value = 0.05478559966737641
$ws.Range("J6").Value = 0.05478559966737641
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 25.53984366666667
$ws.Range("N6").Value = 76.619531
$ws.Range("O6").Value = 0.5588017999566871
$ws.Range("P6").Value = 0.5588017999566871
$ws.Range("Q6").Value = 98.25728667861745
$ws.Range("R6").Value = 884.315580107557
$ws.Range("S6").Value = 0.03061429170583642
$ws.Range("T6").Value = 0.03061429170583642

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.74214666666667
$ws.Range("H7").Value = 53.22644
$ws.Range("I7").Value = 0.2526539265634818
$ws.Range("J7").Value = 0.2526539265634818
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.407369
$ws.Range("N7").Value = 7.222107
$ws.Range("O7").Value = 0.05267229306167105
$ws.Range("P7").Value = 0.05267229306167105
$ws.Range("Q7").Value = 42.71189387878667
$ws.Range("R7").Value = 384.40704490908
$ws.Range("S7").Value = 0.01330786166313363
$ws.Range("T7").Value = 0.01330786166313363

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 17.74214666666667
$ws.Range("H8").Value = 53.22644
$ws.Range("I8").Value = 0.2526539265634818
$ws.Range("J8").Value = 0.2526539265634818
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.67754066666667
$ws.Range("N8").Value = 35.032622
$ws.Range("O8").Value = 0.2555000268900398
$ws.Range("P8").Value = 0.2555000268900398
$ws.Range("Q8").Value = 207.1846392139645
$ws.Range("R8").Value = 1864.66175292568
$ws.Range("S8").Value = 0.06455308503084373
$ws.Range("T8").Value = 0.06455308503084374

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 17.74214666666667
$ws.Range("H9").Value = 53.22644
$ws.Range("I9").Value = 0.2526539265634818
$ws.Range("J9").Value = 0.2526539265634818
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.130105
$ws.Range("N9").Value = 9.390315
$ws.Range("O9").Value = 0.06848547433891598
$ws.Range("P9").Value = 0.06848547433891598
$ws.Range("Q9").Value = 55.53478199206666
$ws.Range("R9").Value = 499.8130379286
$ws.Range("S9").Value = 0.0173031240042897
$ws.Range("T9").Value = 0.0173031240042897

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 17.74214666666667
$ws.Range("H10").Value = 53.22644
$ws.Range("I10").Value = 0.2526539265634818
$ws.Range("J10").Value = 0.2526539265634818
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.949797
$ws.Range("N10").Value = 8.849391
$ws.Range("O10").Value = 0.06454040575268606
$ws.Range("P10").Value = 0.06454040575268606
$ws.Range("Q10").Value = 52.33573101089334
$ws.Range("R10").Value = 471.02157909804
$ws.Range("S10").Value = 0.01630638693541646
$ws.Range("T10").Value = 0.01630638693541646

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 17.74214666666667
$ws.Range("H11").Value = 53.22644
$ws.Range("I11").Value = 0.2526539265634818
$ws.Range("J11").Value = 0.2526539265634818
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 25.53984366666667
$ws.Range("N11").Value = 76.619531
$ws.Range("O11").Value = 0.5588017999566871
$ws.Range("P11").Value = 0.5588017999566871
$ws.Range("Q11").Value = 453.1316521777377
$ws.Range("R11").Value = 4078.18486959964
$ws.Range("S11").Value = 0.1411834689297983
$ws.Range("T11").Value = 0.1411834689297983

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 21.79134566666667
$ws.Range("H12").Value = 65.374037
$ws.Range("I12").Value = 0.3103158344491261
$ws.Range("J12").Value = 0.3103158344491261
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.407369
$ws.Range("N12").Value = 7.222107
$ws.Range("O12").Value = 0.05267229306167105
$ws.Range("P12").Value = 0.05267229306167105
$ws.Range("Q12").Value = 52.45981002621767
$ws.Range("R12").Value = 472.138290235959
$ws.Range("S12").Value = 0.01634504657378136
$ws.Range("T12").Value = 0.01634504657378136

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 21.79134566666667
$ws.Range("H13").Value = 65.374037
$ws.Range("I13").Value = 0.3103158344491261
$ws.Range("J13").Value = 0.3103158344491261
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 11.67754066666667
$ws.Range("N13").Value = 35.032622
$ws.Range("O13").Value = 0.2555000268900398
$ws.Range("P13").Value = 0.2555000268900398
$ws.Range("Q13").Value = 254.4693252038905
$ws.Range("R13").Value = 2290.223926835014
$ws.Range("S13").Value = 0.07928570404615684
$ws.Range("T13").Value = 0.07928570404615685

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 21.79134566666667
$ws.Range("H14").Value = 65.374037
$ws.Range("I14").Value = 0.3103158344491261
$ws.Range("J14").Value = 0.3103158344491261
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.130105
$ws.Range("N14").Value = 9.390315
$ws.Range("O14").Value = 0.06848547433891598
$ws.Range("P14").Value = 0.06848547433891598
$ws.Range("Q14").Value = 68.20920002796167
$ws.Range("R14").Value = 613.882800251655
$ws.Range("S14").Value = 0.02125212711712492
$ws.Range("T14").Value = 0.02125212711712492

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 21.79134566666667
$ws.Range("H15").Value = 65.374037
$ws.Range("I15").Value = 0.3103158344491261
$ws.Range("J15").Value = 0.3103158344491261
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.949797
$ws.Range("N15").Value = 8.849391
$ws.Range("O15").Value = 0.06454040575268606
$ws.Range("P15").Value = 0.06454040575268606
$ws.Range("Q15").Value = 64.28004607349634
$ws.Range("R15").Value = 578.520414661467
$ws.Range("S15").Value = 0.02002790986682995
$ws.Range("T15").Value = 0.02002790986682995

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 21.79134566666667
$ws.Range("H16").Value = 65.374037
$ws.Range("I16").Value = 0.3103158344491261
$ws.Range("J16").Value = 0.3103158344491261
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 25.53984366666667
$ws.Range("N16").Value = 76.619531
$ws.Range("O16").Value = 0.5588017999566871
$ws.Range("P16").Value = 0.5588017999566871
$ws.Range("Q16").Value = 556.5475616129609
$ws.Range("R16").Value = 5008.928054516647
$ws.Range("S16").Value = 0.173405046845233
$ws.Range("T16").Value = 0.173405046845233

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 20.94549
$ws.Range("H17").Value = 62.83647
$ws.Range("I17").Value = 0.2982705752420869
$ws.Range("J17").Value = 0.2982705752420869
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.407369
$ws.Range("N17").Value = 7.222107
$ws.Range("O17").Value = 0.05267229306167105
$ws.Range("P17").Value = 0.05267229306167105
$ws.Range("Q17").Value = 50.42352331581
$ws.Range("R17").Value = 453.81170984229
$ws.Range("S17").Value = 0.01571059515082441
$ws.Range("T17").Value = 0.01571059515082441

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 20.94549
$ws.Range("H18").Value = 62.83647
$ws.Range("I18").Value = 0.2982705752420869
$ws.Range("J18").Value = 0.2982705752420869
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 11.67754066666667
$ws.Range("N18").Value = 35.032622
$ws.Range("O18").Value = 0.2555000268900398
$ws.Range("P18").Value = 0.2555000268900398
$ws.Range("Q18").Value = 244.59181125826
$ws.Range("R18").Value = 2201.32630132434
$ws.Range("S18").Value = 0.07620813999486084
$ws.Range("T18").Value = 0.07620813999486085

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 20.94549
$ws.Range("H19").Value = 62.83647
$ws.Range("I19").Value = 0.2982705752420869
$ws.Range("J19").Value = 0.2982705752420869
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 3.130105
$ws.Range("N19").Value = 9.390315
$ws.Range("O19").Value = 0.06848547433891598
$ws.Range("P19").Value = 0.06848547433891598
$ws.Range("Q19").Value = 65.56158297645
$ws.Range("R19").Value = 590.0542467880499
$ws.Range("S19").Value = 0.02042720182679566
$ws.Range("T19").Value = 0.02042720182679566

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 20.94549
$ws.Range("H20").Value = 62.83647
$ws.Range("I20").Value = 0.2982705752420869
$ws.Range("J20").Value = 0.2982705752420869
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 2.949797
$ws.Range("N20").Value = 8.849391
$ws.Range("O20").Value = 0.06454040575268606
$ws.Range("P20").Value = 0.06454040575268606
$ws.Range("Q20").Value = 61.78494356553
$ws.Range("R20").Value = 556.06449208977
$ws.Range("S20").Value = 0.01925050395021137
$ws.Range("T20").Value = 0.01925050395021137

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 20.94549
$ws.Range("H21").Value = 62.83647
$ws.Range("I21").Value = 0.2982705752420869
$ws.Range("J21").Value = 0.2982705752420869
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 25.53984366666667
$ws.Range("N21").Value = 76.619531
$ws.Range("O21").Value = 0.5588017999566871
$ws.Range("P21").Value = 0.5588017999566871
$ws.Range("Q21").Value = 534.94454012173
$ws.Range("R21").Value = 4814.500861095569
$ws.Range("S21").Value = 0.1666741343193947
$ws.Range("T21").Value = 0.1666741343193947

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 5.896920666666666
$ws.Range("H22").Value = 17.690762
$ws.Range("I22").Value = 0.08397406407792883
$ws.Range("J22").Value = 0.08397406407792883
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 2.407369
$ws.Range("N22").Value = 7.222107
$ws.Range("O22").Value = 0.05267229306167105
$ws.Range("P22").Value = 0.05267229306167105
$ws.Range("Q22").Value = 14.19606400839267
$ws.Range("R22").Value = 127.764576075534
$ws.Range("S22").Value = 0.00442310651269221
$ws.Range("T22").Value = 0.00442310651269221

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 5.896920666666666
$ws.Range("H23").Value = 17.690762
$ws.Range("I23").Value = 0.08397406407792883
$ws.Range("J23").Value = 0.08397406407792883
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 11.67754066666667
$ws.Range("N23").Value = 35.032622
$ws.Range("O23").Value = 0.2555000268900398
$ws.Range("P23").Value = 0.2555000268900398
$ws.Range("Q23").Value = 68.86153089310712
$ws.Range("R23").Value = 619.753778037964
$ws.Range("S23").Value = 0.02145537562997674
$ws.Range("T23").Value = 0.02145537562997674

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 5.896920666666666
$ws.Range("H24").Value = 17.690762
$ws.Range("I24").Value = 0.08397406407792883
$ws.Range("J24").Value = 0.08397406407792883
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 3.130105
$ws.Range("N24").Value = 9.390315
$ws.Range("O24").Value = 0.06848547433891598
$ws.Range("P24").Value = 0.06848547433891598
$ws.Range("Q24").Value = 18.45798086333667
$ws.Range("R24").Value = 166.12182777003
$ws.Range("S24").Value = 0.005751003610543481
$ws.Range("T24").Value = 0.005751003610543481

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 5.896920666666666
$ws.Range("H25").Value = 17.690762
$ws.Range("I25").Value = 0.08397406407792883
$ws.Range("J25").Value = 0.08397406407792883
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 2.949797
$ws.Range("N25").Value = 8.849391
$ws.Range("O25").Value = 0.06454040575268606
$ws.Range("P25").Value = 0.06454040575268606
$ws.Range("Q25").Value = 17.39471889177133
$ws.Range("R25").Value = 156.552470025942
$ws.Range("S25").Value = 0.005419720168291586
$ws.Range("T25").Value = 0.005419720168291586

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 5.896920666666666
$ws.Range("H26").Value = 17.690762
$ws.Range("I26").Value = 0.08397406407792883
$ws.Range("J26").Value = 0.08397406407792883
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 25.53984366666667
$ws.Range("N26").Value = 76.619531
$ws.Range("O26").Value = 0.5588017999566871
$ws.Range("P26").Value = 0.5588017999566871
$ws.Range("Q26").Value = 150.6064319414024
$ws.Range("R26").Value = 1355.457887472622
$ws.Range("S26").Value = 0.04692485815642481
$ws.Range("T26").Value = 0.04692485815642481

